$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price cells (column D) that are about to receive new,
# numeric-looking text (e.g. "1.00", "522.54") to Text format first.
# Without this, Excel auto-converts such strings to numbers (losing
# trailing zeros / thousands-dot formatting), whereas the workbook
# stores these as plain text cells.
$ws.Range('D2:D6').NumberFormat = "@"
$ws.Range('D8:D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14:D16').NumberFormat = "@"
$ws.Range('D18:D22').NumberFormat = "@"
$ws.Range('D24:D29').NumberFormat = "@"
$ws.Range('D31:D36').NumberFormat = "@"
$ws.Range('D38:D48').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '58.676.65'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '2.486.41'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '522.54'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').Value = '133.01'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.559'
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('D9').Value = '2.521.51'
$ws.Range('E9').Value = '  +2.43%  '
$ws.Range('D10').Value = '0.0976'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '5.17'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('D14').Value = '2.957.03'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').Value = '58.396.79'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').Value = '22.16'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '2.505.34'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').Value = '10.67'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('D20').Value = '322.55'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').Value = '6.12'
$ws.Range('E22').Value = '  +7.03%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '64.43'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').Value = '0.402'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = '7.40'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').Value = '0.0₃0756'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('E30').Value = '  +2.23%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '167.86'
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = '6.27'
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('D36').Value = '18.11'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('D38').Value = '3.98'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').Value = '1.48'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').Value = '36.14'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '0.781'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('D42').Value = '3.51'
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('D43').Value = '278.98'
$ws.Range('E43').Value = '  +2.93%  '
$ws.Range('D44').Value = '5.01'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').Value = '0.599'
$ws.Range('E45').Value = '  +2.14%  '
$ws.Range('D46').Value = '123.65'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('D47').Value = '0.0919'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').Value = '0.0502'
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('E49').Value = '  +1.85%  '
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').Value = '16.97'
$ws.Range('E51').Value = '  +1.03%  '
